# Add a new "120 ohm resistor" line item to the mBom worksheet (row 44),
# mirroring the existing resistor/capacitor rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B44").Value = "120 ohm resistor"
$ws.Range("C44").Value = 402
$ws.Range("D44").Value = "120ohm"
$ws.Range("E44").Value = "1/8 W"
$ws.Range("F44").Value = "RES SMD 120 OHM 5% 1/8W 0402"
$ws.Range("G44").Value = "Vishay Dale"
$ws.Range("H44").Value = "CRCW0402120RJNEDHP"
$ws.Range("I44").Value = "541-120YACT-ND"

$ws.Range("K44").Value = 0.17
$ws.Range("L44").Value = 0.14
$ws.Range("M44").Value = 0.056
$ws.Range("N44").Value = 0.02356

# Move the view/selection down to the newly added row, same as the
# author scrolling to and landing on the new line after typing it in.
[void]$ws.Range("N45").Select()
